# Update "column views / visits" numbers (column F) on the "展览" and
# "全部类型" worksheets, per the gh-pages data refresh at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (row -> new F value)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 136
$ws1.Range("F5").Value  = 6733
$ws1.Range("F6").Value  = 86
$ws1.Range("F7").Value  = 433
$ws1.Range("F8").Value  = 140
$ws1.Range("F9").Value  = 6253
$ws1.Range("F12").Value = 1258
$ws1.Range("F14").Value = 100
$ws1.Range("F18").Value = 368
$ws1.Range("F21").Value = 4581
$ws1.Range("F22").Value = 62
$ws1.Range("F23").Value = 40
$ws1.Range("F24").Value = 33
$ws1.Range("F25").Value = 193
$ws1.Range("F26").Value = 70

# Sheet "全部类型" (row -> new F value)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 136
$ws4.Range("F5").Value  = 6733
$ws4.Range("F6").Value  = 86
$ws4.Range("F7").Value  = 433
$ws4.Range("F8").Value  = 140
$ws4.Range("F9").Value  = 6253
$ws4.Range("F12").Value = 1258
$ws4.Range("F14").Value = 100
$ws4.Range("F18").Value = 368
$ws4.Range("F21").Value = 4581
$ws4.Range("F23").Value = 62
$ws4.Range("F24").Value = 40
$ws4.Range("F25").Value = 33
$ws4.Range("F26").Value = 193
$ws4.Range("F27").Value = 70

$wb.Save()
